# Adds a new "2022-Q3" sheet (with fund-holding detail data) right after
# the "总计" (summary) sheet, shifting all the quarterly sheets that used
# to follow it down by one position, and inserts a matching summary row
# into the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q3.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 23
$total.Range("D2").Value = 5.84

# Give the new A2 the same "index column" style as the rows below it,
# and make sure B2:D2 have no special style (matches the other rows).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("B2:D2").ClearFormats()

# Column A is just a 0-based running row index — re-number it so it
# stays sequential (0..7) after the insert shifted the old rows down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right before "2022-Q2" — this
#    pushes 2022-Q2 .. 2020-Q4 one slot later, matching the diff.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "2022-Q3"

# Copy the header / index-column formatting from the sheet we just
# pushed down, so the new sheet matches the look of its siblings.
$refSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$refSheet.Range("A2").Copy()
$ws.Range("A2:A24").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Column B (fund codes, which may have leading zeros) and columns D:F
# (and most of G) hold numeric-looking values that must be stored as
# text, like the source data — pre-format them as Text before assigning
# so they aren't coerced into numbers.
$ws.Range("B2:B24").NumberFormat = "@"
$ws.Range("D2:F24").NumberFormat = "@"
$ws.Range("G2:G22").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "206009"
$ws.Range("C2").Value = "鹏华新兴产业混合"
$ws.Range("D2").Value = "41.17"
$ws.Range("E2").Value = "89.21"
$ws.Range("F2").Value = "6.19"
$ws.Range("G2").Value = "2.5484"
$ws.Range("H2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "010549"
$ws.Range("C3").Value = "富国均衡策略混合"
$ws.Range("D3").Value = "23.14"
$ws.Range("E3").Value = "86.21"
$ws.Range("F3").Value = "2.90"
$ws.Range("G3").Value = "0.6711"
$ws.Range("H3").Value = 6

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "005760"
$ws.Range("C4").Value = "富国周期优势混合A"
$ws.Range("D4").Value = "22.68"
$ws.Range("E4").Value = "86.26"
$ws.Range("F4").Value = "2.77"
$ws.Range("G4").Value = "0.6282"
$ws.Range("H4").Value = 9

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "002851"
$ws.Range("C5").Value = "南方品质优选灵活配置混合A"
$ws.Range("D5").Value = "11.81"
$ws.Range("E5").Value = "71.72"
$ws.Range("F5").Value = "4.30"
$ws.Range("G5").Value = "0.5078"
$ws.Range("H5").Value = 6

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "159996"
$ws.Range("C6").Value = "国泰中证全指家用电器ETF"
$ws.Range("D6").Value = "21.72"
$ws.Range("E6").Value = "97.52"
$ws.Range("F6").Value = "2.16"
$ws.Range("G6").Value = "0.4692"
$ws.Range("H6").Value = 8

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "000471"
$ws.Range("C7").Value = "富国城镇发展股票"
$ws.Range("D7").Value = "9.89"
$ws.Range("E7").Value = "85.67"
$ws.Range("F7").Value = "2.70"
$ws.Range("G7").Value = "0.2670"
$ws.Range("H7").Value = 7

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "560880"
$ws.Range("C8").Value = "广发中证全指家用电器ETF"
$ws.Range("D8").Value = "9.25"
$ws.Range("E8").Value = "99.30"
$ws.Range("F8").Value = "2.18"
$ws.Range("G8").Value = "0.2016"
$ws.Range("H8").Value = 8

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "011956"
$ws.Range("C9").Value = "鹏华新能源精选混合A"
$ws.Range("D9").Value = "7.18"
$ws.Range("E9").Value = "85.90"
$ws.Range("F9").Value = "2.72"
$ws.Range("G9").Value = "0.1953"
$ws.Range("H9").Value = 10

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "011957"
$ws.Range("C10").Value = "鹏华新能源精选混合C"
$ws.Range("D10").Value = "4.87"
$ws.Range("E10").Value = "85.90"
$ws.Range("F10").Value = "2.72"
$ws.Range("G10").Value = "0.1325"
$ws.Range("H10").Value = 10

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "673081"
$ws.Range("C11").Value = "西部利得祥运灵活配置混合A"
$ws.Range("D11").Value = "3.90"
$ws.Range("E11").Value = "43.71"
$ws.Range("F11").Value = "1.90"
$ws.Range("G11").Value = "0.0741"
$ws.Range("H11").Value = 7

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "001983"
$ws.Range("C12").Value = "中邮低碳经济灵活配置混合"
$ws.Range("D12").Value = "0.50"
$ws.Range("E12").Value = "91.21"
$ws.Range("F12").Value = "6.99"
$ws.Range("G12").Value = "0.0350"
$ws.Range("H12").Value = 5

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "561120"
$ws.Range("C13").Value = "富国中证全指家用电器ETF"
$ws.Range("D13").Value = "1.27"
$ws.Range("E13").Value = "99.27"
$ws.Range("F13").Value = "2.18"
$ws.Range("G13").Value = "0.0277"
$ws.Range("H13").Value = 8

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "010447"
$ws.Range("C14").Value = "中邮未来成长混合A"
$ws.Range("D14").Value = "0.43"
$ws.Range("E14").Value = "91.79"
$ws.Range("F14").Value = "5.01"
$ws.Range("G14").Value = "0.0215"
$ws.Range("H14").Value = 3

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "015005"
$ws.Range("C15").Value = "中邮能源革新混合C"
$ws.Range("D15").Value = "0.40"
$ws.Range("E15").Value = "91.38"
$ws.Range("F15").Value = "4.70"
$ws.Range("G15").Value = "0.0188"
$ws.Range("H15").Value = 8

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "011565"
$ws.Range("C16").Value = "富国周期优势混合C"
$ws.Range("D16").Value = "0.41"
$ws.Range("E16").Value = "86.26"
$ws.Range("F16").Value = "2.77"
$ws.Range("G16").Value = "0.0114"
$ws.Range("H16").Value = 9

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "010765"
$ws.Range("C17").Value = "国寿安保华丰混合A"
$ws.Range("D17").Value = "0.45"
$ws.Range("E17").Value = "87.38"
$ws.Range("F17").Value = "2.12"
$ws.Range("G17").Value = "0.0095"
$ws.Range("H17").Value = 10

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "015004"
$ws.Range("C18").Value = "中邮能源革新混合A"
$ws.Range("D18").Value = "0.12"
$ws.Range("E18").Value = "91.38"
$ws.Range("F18").Value = "4.70"
$ws.Range("G18").Value = "0.0056"
$ws.Range("H18").Value = 8

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "165524"
$ws.Range("C19").Value = "信诚中证智能家居指数（LOF）A"
$ws.Range("D19").Value = "0.35"
$ws.Range("E19").Value = "91.70"
$ws.Range("F19").Value = "1.40"
$ws.Range("G19").Value = "0.0049"
$ws.Range("H19").Value = 3

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "010448"
$ws.Range("C20").Value = "中邮未来成长混合C"
$ws.Range("D20").Value = "0.06"
$ws.Range("E20").Value = "91.79"
$ws.Range("F20").Value = "5.01"
$ws.Range("G20").Value = "0.0030"
$ws.Range("H20").Value = 3

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "013084"
$ws.Range("C21").Value = "信诚中证智能家居指数（LOF）C"
$ws.Range("D21").Value = "0.15"
$ws.Range("E21").Value = "91.70"
$ws.Range("F21").Value = "1.40"
$ws.Range("G21").Value = "0.0021"
$ws.Range("H21").Value = 3

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "673083"
$ws.Range("C22").Value = "西部利得祥运灵活配置混合C"
$ws.Range("D22").Value = "0.04"
$ws.Range("E22").Value = "43.71"
$ws.Range("F22").Value = "1.90"
$ws.Range("G22").Value = "0.0008"
$ws.Range("H22").Value = 7

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "010766"
$ws.Range("C23").Value = "国寿安保华丰混合C"
$ws.Range("D23").Value = "0.00"
$ws.Range("E23").Value = "87.38"
$ws.Range("F23").Value = "2.12"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 10

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "013501"
$ws.Range("C24").Value = "南方品质优选灵活配置混合C"
$ws.Range("D24").Value = "0.00"
$ws.Range("E24").Value = "71.72"
$ws.Range("F24").Value = "4.30"
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 6
